# Remove the trailing "PS" section (the "PS" heading paragraph, the
# "Upon reflection, I'm pretty sure Mitchell O'Hara-Wild ..." paragraph
# with its hyperlink, and the final empty paragraph) that used to sit
# just before the section properties at the end of the document body.

$d = $word.ActiveDocument

# Locate the standalone "PS" paragraph near the end of the document.
$psIndex = -1
$total = $d.Paragraphs.Count
for ($i = $total; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    $t = $t.Trim()
    if ($t -eq "PS") {
        $psIndex = $i
        break
    }
}

if ($psIndex -ge 1) {
    $startRange = $d.Paragraphs.Item($psIndex).Range.Start

    # Delete everything from the start of the "PS" paragraph through to
    # the end of the document's main story (removes the "PS" paragraph,
    # the following "Upon reflection..." paragraph with its hyperlink,
    # and the trailing empty paragraph), leaving the section break intact.
    $endRange = $d.Content.End

    $r = $d.Range($startRange, $endRange)
    $r.Delete()
}
